$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header labels "1n" and "2n" above columns K and L, right-aligned
$ws.Range("K1").Value = "1n"
$ws.Range("L1").Value = "2n"
$ws.Range("K1:L1").HorizontalAlignment = -4152

# DAC clock frequency increase: 3760 -> 4160
$ws.Range("C2").Formula = "=48000*4160/1000000"

# Add TXCO correction (2 nyquist) formula in L2
$ws.Range("L2").Formula = "=K2*2"
